# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures across the Leve Profits sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)
# with the latest market board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1800.2307
$ws.Range("I86").Value = 1800.3334
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 1800.3334
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -677.3334
$ws.Range("N86").Value = -4046

$ws.Range("H89").Value = 1800.2307
$ws.Range("I89").Value = 1800.3334
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 9001.666999999999
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -3385.666999999999
$ws.Range("N89").Value = -20232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 1236.1428
$ws.Range("I31").Value = 1236.1428
$ws.Range("K31").Value = 1236.1428
$ws.Range("M31").Value = -942.1428000000001

$ws.Range("H80").Value = 27800
$ws.Range("J80").Value = 27800
$ws.Range("L80").Value = 27800
$ws.Range("N80").Value = -29796

$ws.Range("H83").Value = 27800
$ws.Range("J83").Value = 27800
$ws.Range("L83").Value = 83400
$ws.Range("N83").Value = -93384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1871.9166
$ws.Range("I16").Value = 1140
$ws.Range("J16").Value = 2394.7144
$ws.Range("K16").Value = 1140
$ws.Range("L16").Value = 2394.7144
$ws.Range("M16").Value = -853
$ws.Range("N16").Value = -2968.7144

$ws.Range("H58").Value = 1716.7273
$ws.Range("I58").Value = 1211.4
$ws.Range("J58").Value = 2799.5715
$ws.Range("K58").Value = 1211.4
$ws.Range("L58").Value = 2799.5715
$ws.Range("M58").Value = -1008.4
$ws.Range("N58").Value = -3205.5715

$ws.Range("H62").Value = 3010.8
$ws.Range("I62").Value = 2939.8
$ws.Range("J62").Value = 3081.8
$ws.Range("K62").Value = 2939.8
$ws.Range("L62").Value = 3081.8
$ws.Range("M62").Value = -2315.8
$ws.Range("N62").Value = -4329.8

$ws.Range("H65").Value = 3010.8
$ws.Range("I65").Value = 2939.8
$ws.Range("J65").Value = 3081.8
$ws.Range("K65").Value = 14699
$ws.Range("L65").Value = 15409
$ws.Range("M65").Value = -11579
$ws.Range("N65").Value = -21649

$ws.Range("H80").Value = 18000
$ws.Range("J80").Value = 18000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -20246

$ws.Range("H83").Value = 18000
$ws.Range("J83").Value = 18000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -65232

$ws.Range("H113").Value = 1871.9166
$ws.Range("I113").Value = 1140
$ws.Range("J113").Value = 2394.7144
$ws.Range("K113").Value = 1140
$ws.Range("L113").Value = 2394.7144
$ws.Range("M113").Value = 1030
$ws.Range("N113").Value = -6734.7144

$ws.Range("H132").Value = 2030.359
$ws.Range("I132").Value = 1529.1538
$ws.Range("J132").Value = 3032.7693
$ws.Range("K132").Value = 4587.4614
$ws.Range("L132").Value = 9098.3079
$ws.Range("M132").Value = -2057.4614
$ws.Range("N132").Value = -14158.3079

$ws.Range("H136").Value = 1716.7273
$ws.Range("I136").Value = 1211.4
$ws.Range("J136").Value = 2799.5715
$ws.Range("K136").Value = 3634.2
$ws.Range("L136").Value = 8398.7145
$ws.Range("M136").Value = -1084.2
$ws.Range("N136").Value = -13498.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2666.8572
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2666.8572
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 8000.571599999999
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -8224.571599999999

$ws.Range("H107").Value = 667079.2
$ws.Range("I107").Value = 152.42857
$ws.Range("J107").Value = 1250640.1
$ws.Range("K107").Value = 457.28571
$ws.Range("L107").Value = 3751920.3
$ws.Range("M107").Value = 1462.71429
$ws.Range("N107").Value = -3755760.3

$ws.Range("H132").Value = 12699197
$ws.Range("I132").Value = 859
$ws.Range("J132").Value = 16162380
$ws.Range("K132").Value = 7731
$ws.Range("L132").Value = 145461420
$ws.Range("M132").Value = -5201
$ws.Range("N132").Value = -145466480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 14900
$ws.Range("J15").Value = 14900
$ws.Range("L15").Value = 14900
$ws.Range("N15").Value = -15476

$ws.Range("H81").Value = 14900
$ws.Range("J81").Value = 14900
$ws.Range("L81").Value = 14900
$ws.Range("N81").Value = -16896

$ws.Range("H84").Value = 14900
$ws.Range("J84").Value = 14900
$ws.Range("L84").Value = 44700
$ws.Range("N84").Value = -54684

$ws.Range("H122").Value = 3575202.5
$ws.Range("I122").Value = 10001440
$ws.Range("J122").Value = 5070.4443
$ws.Range("K122").Value = 30004320
$ws.Range("L122").Value = 15211.3329
$ws.Range("M122").Value = -30001870
$ws.Range("N122").Value = -20111.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 8000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -8224

$ws.Range("H16").Value = 1040
$ws.Range("I16").Value = 1126.4445
$ws.Range("J16").Value = 651
$ws.Range("K16").Value = 1126.4445
$ws.Range("L16").Value = 651
$ws.Range("M16").Value = -956.4445000000001
$ws.Range("N16").Value = -991

$ws.Range("H46").Value = 742
$ws.Range("I46").Value = 598
$ws.Range("J46").Value = 850
$ws.Range("K46").Value = 598
$ws.Range("L46").Value = 850
$ws.Range("M46").Value = -410
$ws.Range("N46").Value = -1226

$ws.Range("H68").Value = 2166.4546
$ws.Range("I68").Value = 2021.2941
$ws.Range("K68").Value = 2021.2941
$ws.Range("M68").Value = -1272.2941

$ws.Range("H71").Value = 2166.4546
$ws.Range("I71").Value = 2021.2941
$ws.Range("K71").Value = 10106.4705
$ws.Range("M71").Value = -6362.470499999999

$ws.Range("H136").Value = 2441.9167
$ws.Range("I136").Value = 1169
$ws.Range("J136").Value = 3460.25
$ws.Range("K136").Value = 3507
$ws.Range("L136").Value = 10380.75
$ws.Range("M136").Value = -957
$ws.Range("N136").Value = -15480.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 18325
$ws.Range("J86").Value = 18325
$ws.Range("L86").Value = 18325
$ws.Range("N86").Value = -20571

$ws.Range("H89").Value = 18325
$ws.Range("J89").Value = 18325
$ws.Range("L89").Value = 91625
$ws.Range("N89").Value = -102857
